# 9th Stab- Cosmetic Changes
# Insert two new weekly columns ("Jun_17" and "Jun_15") before the existing
# "this week" column, shifting prior weeks (Jun_13, Jun_10) to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new columns at C:D (old column C and its data moves to column E).
$ws.Columns("C:D").Insert()

# Match the cosmetic column width formatting used by column C (width 8.0,
# customWidth true) on the two freshly inserted columns.
$ws.Columns("C:D").ColumnWidth = 7.166666666666667

# New header row: B1 becomes the newest week, C1/D1 are the two newly
# inserted weeks, and the old header (now in D1 after the insert shift)
# keeps its original value.
$oldB1 = $ws.Range("B1").Value2
$ws.Range("D1").Value2 = $oldB1
$ws.Range("B1").Value2 = "Jun_17"
$ws.Range("C1").Value2 = "Jun_15"

# Fill the two new weekly columns with the same default rating ("UN") used
# in column B for every analyst row.
for ($r = 2; $r -le 27; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value2 = $bVal
    $ws.Cells.Item($r, 4).Value2 = $bVal
}
